$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new/updated data for rows 84-94 (columns B "Morning" and C "Evening")
$ws.Range("C84").Value = 116.4

$ws.Range("B85").Value = 115.4
$ws.Range("C85").Value = 116.4

$ws.Range("B86").Value = 115.4
$ws.Range("C86").Value = 116.4

$ws.Range("B87").Value = 114.9
$ws.Range("C87").Value = 116.4

$ws.Range("B88").Value = 115.4
$ws.Range("C88").Value = 117

$ws.Range("B89").Value = 116.1
$ws.Range("C89").Value = 116.4

$ws.Range("B90").Value = 115.9
$ws.Range("C90").Value = 116

$ws.Range("B91").Value = 114.7
$ws.Range("C91").Value = 116

$ws.Range("B92").Value = 114.5
$ws.Range("C92").Value = 116.2

$ws.Range("B93").Value = 114.6
$ws.Range("C93").Value = 116.9

$ws.Range("B94").Value = 115.7

# Update the view to reflect the new scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 75
$ws.Range("B94").Select()
